# Scheduled market-price refresh: updates the currentAveragePrice* /
# LevePrice* / LeveProfit* columns (H-N) for a handful of leve rows across
# the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets with freshly pulled numbers.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 5848.7
$ws.Range("I51").Value = 5997
$ws.Range("K51").Value = 5997
$ws.Range("M51").Value = -5513

$ws.Range("H88").Value = 20870102
$ws.Range("I88").Value = 111113816
$ws.Range("J88").Value = 44630.31
$ws.Range("K88").Value = 111113816
$ws.Range("L88").Value = 44630.31
$ws.Range("M88").Value = -111113410
$ws.Range("N88").Value = -45442.31

$ws.Range("H91").Value = 20870102
$ws.Range("I91").Value = 111113816
$ws.Range("J91").Value = 44630.31
$ws.Range("K91").Value = 111113816
$ws.Range("L91").Value = 44630.31
$ws.Range("M91").Value = -111112412
$ws.Range("N91").Value = -47438.31

$ws.Range("H100").Value = 1646.9333
$ws.Range("I100").Value = 2911.75
$ws.Range("J100").Value = 201.42857
$ws.Range("K100").Value = 2911.75
$ws.Range("L100").Value = 201.42857
$ws.Range("M100").Value = -2370.75
$ws.Range("N100").Value = -1283.42857

$ws.Range("H103").Value = 867.41174
$ws.Range("I103").Value = 456.83334
$ws.Range("J103").Value = 1091.3636
$ws.Range("K103").Value = 1370.50002
$ws.Range("L103").Value = 3274.0908
$ws.Range("M103").Value = -784.5000199999999
$ws.Range("N103").Value = -4446.0908

$ws.Range("H121").Value = 988.75
$ws.Range("J121").Value = 1168.3334
$ws.Range("L121").Value = 3505.0002
$ws.Range("N121").Value = -6999.0002

$ws.Range("H132").Value = 1091.9767
$ws.Range("I132").Value = 1082.8049
$ws.Range("J132").Value = 1280
$ws.Range("K132").Value = 3248.4147
$ws.Range("L132").Value = 3840
$ws.Range("M132").Value = -718.4147000000003
$ws.Range("N132").Value = -8900

$ws.Range("H137").Value = 6012.3667
$ws.Range("I137").Value = 2863.353
$ws.Range("K137").Value = 8590.059000000001
$ws.Range("M137").Value = -6040.059000000001

$ws.Range("H138").Value = 1055702.4
$ws.Range("I138").Value = 1254.86
$ws.Range("J138").Value = 2227310.8
$ws.Range("K138").Value = 3764.58
$ws.Range("L138").Value = 6681932.399999999
$ws.Range("M138").Value = 1375.42
$ws.Range("N138").Value = -6692212.399999999

$ws.Range("H141").Value = 4304.4346
$ws.Range("I141").Value = 3857.238
$ws.Range("J141").Value = 9000
$ws.Range("K141").Value = 11571.714
$ws.Range("L141").Value = 27000
$ws.Range("M141").Value = -6391.714
$ws.Range("N141").Value = -37360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2863591
$ws.Range("I32").Value = 3035452
$ws.Range("K32").Value = 3035452
$ws.Range("M32").Value = -3035165

$ws.Range("H45").Value = 4788.5454
$ws.Range("I45").Value = 1910.75
$ws.Range("K45").Value = 1910.75
$ws.Range("M45").Value = -1533.75

$ws.Range("H61").Value = 43485776
$ws.Range("I61").Value = 3523.9333
$ws.Range("J61").Value = 125015000
$ws.Range("K61").Value = 3523.9333
$ws.Range("L61").Value = 125015000
$ws.Range("M61").Value = -3311.9333
$ws.Range("N61").Value = -125015424

$ws.Range("H74").Value = 39448.035
$ws.Range("I74").Value = 65165.812
$ws.Range("K74").Value = 65165.812
$ws.Range("M74").Value = -64291.812

$ws.Range("H77").Value = 39448.035
$ws.Range("I77").Value = 65165.812
$ws.Range("K77").Value = 325829.06
$ws.Range("M77").Value = -321461.06

$ws.Range("H132").Value = 4736.3726
$ws.Range("I132").Value = 2264.6897
$ws.Range("J132").Value = 7994.5
$ws.Range("K132").Value = 6794.0691
$ws.Range("L132").Value = 23983.5
$ws.Range("M132").Value = -4264.0691
$ws.Range("N132").Value = -29043.5

$ws.Range("H136").Value = 43485776
$ws.Range("I136").Value = 3523.9333
$ws.Range("J136").Value = 125015000
$ws.Range("K136").Value = 10571.7999
$ws.Range("L136").Value = 375045000
$ws.Range("M136").Value = -8021.7999
$ws.Range("N136").Value = -375050100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 11118936
$ws.Range("I20").Value = 16673316
$ws.Range("K20").Value = 16673316
$ws.Range("M20").Value = -16673069

$ws.Range("H21").Value = 18613.75
$ws.Range("J21").Value = 18613.75
$ws.Range("L21").Value = 18613.75
$ws.Range("N21").Value = -19085.75

$ws.Range("H27").Value = 80000
$ws.Range("J27").Value = 80000
$ws.Range("L27").Value = 80000
$ws.Range("N27").Value = -80384

$ws.Range("H94").Value = 2554.2104
$ws.Range("I94").Value = 1486.1428
$ws.Range("K94").Value = 1486.1428
$ws.Range("M94").Value = -1035.1428

$ws.Range("H105").Value = 2390.4092
$ws.Range("I105").Value = 2002
$ws.Range("K105").Value = 2002
$ws.Range("M105").Value = -255

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 7144145
$ws.Range("I105").Value = 8929307
$ws.Range("K105").Value = 8929307
$ws.Range("M105").Value = -8927560

$ws.Range("H106").Value = 5000
$ws.Range("J106").Value = 5000
$ws.Range("L106").Value = 5000
$ws.Range("N106").Value = -7524

$ws.Range("H132").Value = 10114.941
$ws.Range("I132").Value = 7392
$ws.Range("K132").Value = 22176
$ws.Range("M132").Value = -19646

$ws.Range("H134").Value = 6195.967
$ws.Range("J134").Value = 9137.412
$ws.Range("L134").Value = 27412.236
$ws.Range("N134").Value = -32482.236

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1547.28
$ws.Range("J5").Value = 3374.375
$ws.Range("L5").Value = 10123.125
$ws.Range("N5").Value = -10347.125

$ws.Range("H98").Value = 1761.1428
$ws.Range("I98").Value = 2453
$ws.Range("J98").Value = 1242.25
$ws.Range("K98").Value = 7359
$ws.Range("L98").Value = 3726.75
$ws.Range("M98").Value = -5861
$ws.Range("N98").Value = -6722.75

$ws.Range("H104").Value = 5337
$ws.Range("I104").Value = 1341.6666
$ws.Range("K104").Value = 4024.9998
$ws.Range("M104").Value = -1403.9998

$ws.Range("H113").Value = 6099.0586
$ws.Range("J113").Value = 10114.556
$ws.Range("L113").Value = 30343.668
$ws.Range("N113").Value = -34683.66800000001

$ws.Range("H125").Value = 5266
$ws.Range("J125").Value = 5266
$ws.Range("L125").Value = 15798
$ws.Range("N125").Value = -25638

$ws.Range("H135").Value = 1547.28
$ws.Range("J135").Value = 3374.375
$ws.Range("L135").Value = 30369.375
$ws.Range("N135").Value = -35439.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 1538668.1
$ws.Range("I2").Value = 58.8
$ws.Range("J2").Value = 2500299
$ws.Range("K2").Value = 58.8
$ws.Range("L2").Value = 2500299
$ws.Range("M2").Value = 54.2
$ws.Range("N2").Value = -2500525

$ws.Range("H43").Value = 1073.8572
$ws.Range("I43").Value = 1073.8572
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 1073.8572
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = -922.8571999999999
$ws.Range("N43").ClearContents()

$ws.Range("H93").Value = 51971.8
$ws.Range("J93").Value = 51971.8
$ws.Range("L93").Value = 51971.8
$ws.Range("N93").Value = -55715.8

$ws.Range("H136").Value = 30636.4
$ws.Range("I136").Value = 25800
$ws.Range("J136").Value = 31603.68
$ws.Range("K136").Value = 77400
$ws.Range("L136").Value = 94811.04000000001
$ws.Range("M136").Value = -74850
$ws.Range("N136").Value = -99911.04000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1569133.1
$ws.Range("I46").Value = 2300203.2
$ws.Range("K46").Value = 2300203.2
$ws.Range("M46").Value = -2300015.2

$ws.Range("H55").Value = 335.08694
$ws.Range("I55").Value = 128.16667
$ws.Range("K55").Value = 128.16667
$ws.Range("M55").Value = 44.83332999999999

$ws.Range("H105").Value = 34666.668
$ws.Range("J105").Value = 34666.668
$ws.Range("L105").Value = 34666.668
$ws.Range("N105").Value = -41654.668

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 4360.357
$ws.Range("I2").Value = 4253.75
$ws.Range("K2").Value = 4253.75
$ws.Range("M2").Value = -4141.75

$ws.Range("H96").Value = 1206
$ws.Range("I96").Value = 1100
$ws.Range("K96").Value = 1100
$ws.Range("M96").Value = 273

